$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated query text for the "Summary" stat query (row 2 header row not affected) ---
$qSummary = @'
SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Gastric Cancer%';
'@

# --- Updated ProjectsTab query (cell B3) ---
$qProjects = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;
'@

# --- Updated GrantsTab query (cell B4) ---
$qGrants = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.project_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@

# --- Updated PublicationsTab query (cell B5) ---
$qPubs = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;
'@

# --- New ProgramsTab query (cell B2), replaces the old one (adds CASE expression + LIKE + lower() sort) ---
$qPrograms = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
    prg.website AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
     CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Gastric Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

# Order matters: the engine appends each newly-written distinct string to the
# end of the shared-strings table (and drops orphaned entries), so writing in
# this order reproduces the target table layout (summary, projects, grants,
# publications, then the new programs query last).
$ws.Range("C2").Value = $qSummary
$ws.Range("B3").Value = $qProjects
$ws.Range("B4").Value = $qGrants
$ws.Range("B5").Value = $qPubs
$ws.Range("B2").Value = $qPrograms

# --- Update the active view/selection (previously B5 was selected with A5 scrolled to top; now C3 is selected with A2 at top) ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C3").Select()
